# quarterly.xlsx update:
#  - drop the oldest quarter column (فصل دوم منتهی به 1399/06) and roll every
#    quarter's data one column to the left (E..M <- F..N)
#  - append the newest quarter (فصل چهارم منتهی به 1401/12) into column N
#  - refresh a couple of the rolled-forward figures that were recomputed by
#    the updated read_price algorithm (columns J on rows 26/27 land on new,
#    not merely shifted, values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header rows (quarter labels), row 8 and row 24 ----
$quarterHeaders = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
$cols = @("E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $quarterHeaders[$i]
    $ws.Range($cols[$i] + "24").Value = $quarterHeaders[$i]
}

# ---- data rows that roll left with the quarter window, with a fresh value appended in N ----
$rowValues = @{
    10 = @(0, 26317, 0, 0, 0, 48508, 14916, 29721, 22318, 20805)
    14 = @(0, 61764, 0, 0, 0, 139139, 12364, 45015, 123214, 148883)
    16 = @(29201, -22100, 57603, 10846, 34719, -37578, 20266, 21745, 21005, 26688)
    17 = @(152360, -19643, 82790, 137505, 96909, 400919, 267499, 244392, 256946, 335870)
    19 = @(137379, 141953, 260345, 144994, 133760, 682486, 609478, -53391, 138136, 929916)
    20 = @(318940, 188291, 400738, 293345, 265388, 1233474, 924523, 287482, 561619, 1462162)
    26 = @(747, 722, 708, 722, 722, 852, 735, 210, 210, 644)
    27 = @(416, 416, 444, 416, 416, 520, 424, 1218, 1218, 788)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
